$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "TestAutomation179"
$ws.Range("B8").Value = "TestAutomation179"
$ws.Range("C8").Value = "Facility_POC179"
$ws.Range("D8").Value = "Facility_POC179"
$ws.Range("E8").Value = "Pharmacy_POC179"
$ws.Range("F8").Value = "Pharmacy_POC179"
$ws.Range("H8").Value = "AlignmentProjectPOC179"

$ws.Range("H16").Select()
